$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 20), copying A19's date formatting/style
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A20").Value = 45776

$ws.Range("B20").Value = "2 hours"
$ws.Range("C20").Value = "update data and plots"
$ws.Range("E20").Value = "N"

$ws.Range("B18").Select()
